$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period-end labels (shift window forward by one quarter) ---
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: publish-date labels (shift window forward by one quarter) ---
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-03-25 (8)"
$ws.Range("F9").Value = "1401-04-28 (2)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-28 (7)"
$ws.Range("J9").Value = "1401-04-28"
$ws.Range("K9").Value = "1401-08-29 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-28"

# --- Row 11: Sales (فروش) ---
$ws.Range("D11").Value = 6720
$ws.Range("E11").Value = 8132
$ws.Range("F11").Value = 11444
$ws.Range("G11").Value = 11379
$ws.Range("H11").Value = 9189
$ws.Range("I11").Value = 8266
$ws.Range("J11").Value = 13311
$ws.Range("K11").Value = 13653
$ws.Range("L11").Value = 11085
$ws.Range("M11").Value = 12316

# --- Row 12: Cost of goods sold ---
$ws.Range("D12").Value = -4151
$ws.Range("E12").Value = -6335
$ws.Range("F12").Value = -8946
$ws.Range("G12").Value = -9042
$ws.Range("H12").Value = -7431
$ws.Range("I12").Value = -7446
$ws.Range("J12").Value = -10251
$ws.Range("K12").Value = -10974
$ws.Range("L12").Value = -8763
$ws.Range("M12").Value = -9968

# --- Row 13: Gross profit ---
$ws.Range("D13").Value = 2569
$ws.Range("E13").Value = 1797
$ws.Range("F13").Value = 2497
$ws.Range("G13").Value = 2338
$ws.Range("H13").Value = 1759
$ws.Range("I13").Value = 819
$ws.Range("J13").Value = 3060
$ws.Range("K13").Value = 2679
$ws.Range("L13").Value = 2322
$ws.Range("M13").Value = 2348

# --- Row 14: Impairment expense ---
$ws.Range("D14").Value = -145
$ws.Range("E14").Value = -348
$ws.Range("F14").Value = -264
$ws.Range("G14").Value = -253
$ws.Range("H14").Value = -269
$ws.Range("I14").Value = -438
$ws.Range("J14").Value = -283
$ws.Range("K14").Value = -364
$ws.Range("L14").Value = -362
$ws.Range("M14").Value = -944

# --- Row 16: Other operating income/expense, net ---
$ws.Range("D16").Value = 646
$ws.Range("E16").Value = -70
$ws.Range("F16").Value = 119
$ws.Range("G16").Value = 349
$ws.Range("H16").Value = 106
$ws.Range("I16").Value = -71
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = -19
$ws.Range("L16").Value = 82
$ws.Range("M16").Value = 312

# --- Row 17: Operating profit ---
$ws.Range("D17").Value = 3070
$ws.Range("E17").Value = 1380
$ws.Range("F17").Value = 2353
$ws.Range("G17").Value = 2433
$ws.Range("H17").Value = 1595
$ws.Range("I17").Value = 311
$ws.Range("J17").Value = 2781
$ws.Range("K17").Value = 2296
$ws.Range("L17").Value = 2042
$ws.Range("M17").Value = 1715

# --- Row 18: Financial expenses ---
$ws.Range("D18").Value = -100
$ws.Range("E18").Value = -187
$ws.Range("F18").Value = -168
$ws.Range("G18").Value = -220
$ws.Range("H18").Value = -249
$ws.Range("I18").Value = -395
$ws.Range("J18").Value = -375
$ws.Range("K18").Value = -388
$ws.Range("L18").Value = -329
$ws.Range("M18").Value = -259

# --- Row 19: Other non-operating income/expense, net ---
$ws.Range("D19").Value = 37
$ws.Range("E19").Value = 62
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 29
$ws.Range("I19").Value = 102
$ws.Range("J19").Value = 5
$ws.Range("K19").Value = 151
$ws.Range("L19").Value = 24
$ws.Range("M19").Value = 206

# --- Row 20: Net profit from continuing ops before tax ---
$ws.Range("D20").Value = 3007
$ws.Range("E20").Value = 1255
$ws.Range("F20").Value = 2258
$ws.Range("G20").Value = 2254
$ws.Range("H20").Value = 1375
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 2412
$ws.Range("K20").Value = 2058
$ws.Range("L20").Value = 1737
$ws.Range("M20").Value = 1662

# --- Row 21: Tax ---
$ws.Range("D21").Value = -369
$ws.Range("E21").Value = 67
$ws.Range("F21").Value = -322
$ws.Range("G21").Value = -515
$ws.Range("H21").Value = -286
$ws.Range("I21").Value = 59
$ws.Range("J21").Value = -423
$ws.Range("K21").Value = 95
$ws.Range("L21").Value = -355
$ws.Range("M21").Value = 211

# --- Row 22: Net profit from continuing ops ---
$ws.Range("D22").Value = 2638
$ws.Range("E22").Value = 1323
$ws.Range("F22").Value = 1936
$ws.Range("G22").Value = 1738
$ws.Range("H22").Value = 1089
$ws.Range("I22").Value = 77
$ws.Range("J22").Value = 1989
$ws.Range("K22").Value = 2153
$ws.Range("L22").Value = 1382
$ws.Range("M22").Value = 1873

# --- Row 24: Net profit ---
$ws.Range("D24").Value = 2638
$ws.Range("E24").Value = 1323
$ws.Range("F24").Value = 1936
$ws.Range("G24").Value = 1738
$ws.Range("H24").Value = 1089
$ws.Range("I24").Value = 77
$ws.Range("J24").Value = 1989
$ws.Range("K24").Value = 2153
$ws.Range("L24").Value = 1382
$ws.Range("M24").Value = 1873

# --- Row 26: Capital ---
$ws.Range("D26").Value = 3312
$ws.Range("E26").Value = 3685
$ws.Range("F26").Value = 3858
$ws.Range("G26").Value = 3448
$ws.Range("H26").Value = 3161
$ws.Range("I26").Value = 3268
$ws.Range("J26").Value = 3063
$ws.Range("K26").Value = 2902
$ws.Range("L26").Value = 5469
$ws.Range("M26").Value = 4180

